# Applies the "456a3b4" data refresh to the 江西-漫展信息 workbook.
# The same set of row-level updates (attendee counts bumped, one listing
# marked cancelled/unsellable, one venue address corrected) must land on
# both the "展览" sheet and the "全部类型" sheet, which mirror each other.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F3: want-to-go count 19 -> 20
    $ws.Range("F3").Value = 20

    # C5: mark the Gao'an event as cancelled
    $ws.Range("C5").Value = "高安·星语动漫嘉年华（取消）"
    # G5: min price 40 -> no longer sellable
    $ws.Range("G5").Value = "不可售"

    # F6: 17 -> 18
    $ws.Range("F6").Value = 18
    # F7: 106 -> 105
    $ws.Range("F7").Value = 105
    # F8: 74 -> 75
    $ws.Range("F8").Value = 75
    # F9: 442 -> 443
    $ws.Range("F9").Value = 443
    # F10: 40 -> 42
    $ws.Range("F10").Value = 42
    # F11: 15 -> 17
    $ws.Range("F11").Value = 17

    # D12: corrected venue address for the Ganzhou event
    $ws.Range("D12").Value = "赞贤路与长征大道交叉口东南200米赣州市少林功夫表演团附近 赣州市体育中心-体育馆"
    # F12: 546 -> 547
    $ws.Range("F12").Value = 547

    # F13: 21 -> 24
    $ws.Range("F13").Value = 24
    # F15: 21 -> 22
    $ws.Range("F15").Value = 22
    # F16: 341 -> 343
    $ws.Range("F16").Value = 343
    # F18: 86 -> 87
    $ws.Range("F18").Value = 87
    # F22: 843 -> 849
    $ws.Range("F22").Value = 849
    # F23: 1366 -> 1369
    $ws.Range("F23").Value = 1369
    # F24: 287 -> 289
    $ws.Range("F24").Value = 289
    # F25: 299 -> 303
    $ws.Range("F25").Value = 303
    # F27: 62 -> 64
    $ws.Range("F27").Value = 64
    # F29: 33 -> 34
    $ws.Range("F29").Value = 34
    # F30: 81 -> 82
    $ws.Range("F30").Value = 82
    # F31: 200 -> 201
    $ws.Range("F31").Value = 201
    # F32: 233 -> 235
    $ws.Range("F32").Value = 235
    # F35: 46 -> 47
    $ws.Range("F35").Value = 47
    # F37: 150 -> 151
    $ws.Range("F37").Value = 151
    # F38: 560 -> 564
    $ws.Range("F38").Value = 564
    # F40: (3379 on 展览 / 3380 on 全部类型) -> 3400 on both
    $ws.Range("F40").Value = 3400
    # F41: 396 -> 399
    $ws.Range("F41").Value = 399
    # F42: 174 -> 176
    $ws.Range("F42").Value = 176
    # F43: 868 -> 875
    $ws.Range("F43").Value = 875
    # F45: 56 -> 57
    $ws.Range("F45").Value = 57
    # F46: 39 -> 40
    $ws.Range("F46").Value = 40
}
